# Update imputed values in Sheet1 to reflect the re-run of the RandomForest
# imputation algorithm (commit: "Update Name of Algo").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.976700000000005
$ws.Range("A9").Value = -22.10050000000002
$ws.Range("D9").Value = -7.722899999999998
$ws.Range("D11").Value = -7.938099999999995
$ws.Range("A18").Value = -22.30910000000001
$ws.Range("A20").Value = -21.56789999999997
$ws.Range("D23").Value = -8.0832
$ws.Range("D24").Value = -7.418799999999997
$ws.Range("D26").Value = -7.549400000000002
$ws.Range("A27").Value = -21.94879999999999
$ws.Range("D34").Value = -7.761100000000002
$ws.Range("A35").Value = -19.64199999999999
$ws.Range("D35").Value = -7.672500000000003
$ws.Range("D48").Value = -7.741399999999999
$ws.Range("D49").Value = -7.971899999999998
$ws.Range("D52").Value = -7.706200000000001
$ws.Range("D66").Value = -7.4034
$ws.Range("D67").Value = -7.526199999999994
$ws.Range("A69").Value = -21.62589999999998
$ws.Range("A76").Value = -19.6443
$ws.Range("A78").Value = -19.09649999999999
$ws.Range("D78").Value = -7.758800000000003
$ws.Range("D80").Value = -8.073000000000002
$ws.Range("A82").Value = -21.67070000000001
$ws.Range("A83").Value = -21.689
$ws.Range("A93").Value = -21.3755
$ws.Range("D99").Value = -7.935099999999998
$ws.Range("D104").Value = -7.737899999999996
